# Generate Report for Archive
#
# The localization status text "Ready for handoff" is now stale - the
# items have moved on to the next pipeline stage, so refresh the status
# everywhere it is reported: the per-locale status columns on the
# "Overview" sheet (columns E and F) and the "Status" column on each of
# the per-locale detail sheets ("zh-cn" and "de-de"). Re-fit the status
# columns afterwards now that the new text is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Columns("E:F").ColumnWidth = 12.5

# --- zh-cn / de-de sheets: Status column (C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2:C3").Value = $newStatus
    $ws.Columns("C:C").ColumnWidth = 12.5
}
